$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 13.5021029201832
$ws.Range("D2").Value = 5.671995549209638
$ws.Range("E2").Value = 13.87812280381266
$ws.Range("F2").Value = 28.95260332316069
$ws.Range("G2").Value = 3.656587908909558
$ws.Range("K2").Value = 16.16283972534575
$ws.Range("L2").Value = 9.264634646378475
$ws.Range("N2").Value = 17.47287962592588
$ws.Range("O2").Value = 25.80717117088754

# Row 3
$ws.Range("C3").Value = 13.40641848683153
$ws.Range("D3").Value = 5.66383323277018
$ws.Range("E3").Value = 13.82006443092378
$ws.Range("F3").Value = 28.94217564587153
$ws.Range("G3").Value = 3.659122450589753
$ws.Range("K3").Value = 15.71195715206656
$ws.Range("L3").Value = 9.268045469743365
$ws.Range("N3").Value = 17.51876110098909
$ws.Range("O3").Value = 25.85783259520933

# Row 4
$ws.Range("C4").Value = 13.35079451305244
$ws.Range("D4").Value = 5.658875181420108
$ws.Range("E4").Value = 13.78758565454162
$ws.Range("F4").Value = 28.94472175422284
$ws.Range("G4").Value = 3.660760630524446
$ws.Range("K4").Value = 15.43065164442201
$ws.Range("L4").Value = 9.271696028239642
$ws.Range("N4").Value = 17.54882500812048
$ws.Range("O4").Value = 25.89592563625709

# Row 5
$ws.Range("C5").Value = 13.32893417067858
$ws.Range("D5").Value = 5.656868821568107
$ws.Range("E5").Value = 13.77515692868987
$ws.Range("F5").Value = 28.94800681638604
$ws.Range("G5").Value = 3.661448881119225
$ws.Range("K5").Value = 15.31506620340165
$ws.Range("L5").Value = 9.273575238775885
$ws.Range("N5").Value = 17.56155301327738
$ws.Range("O5").Value = 25.91319918704575

# Row 6
$ws.Range("C6").Value = 13.32535358554466
$ws.Range("D6").Value = 5.656536535221745
$ws.Range("E6").Value = 13.77314212945512
$ws.Range("F6").Value = 28.94868792564225
$ws.Range("G6").Value = 3.661564415599606
$ws.Range("K6").Value = 15.29582107461559
$ws.Range("L6").Value = 9.273910934922609
$ws.Range("N6").Value = 17.56369530641227
$ws.Range("O6").Value = 25.91617295956658

# Row 7
$ws.Range("C7").Value = 13.35049640361146
$ws.Range("D7").Value = 5.658848065142763
$ws.Range("E7").Value = 13.78741475833352
$ws.Range("F7").Value = 28.94475696273543
$ws.Range("G7").Value = 3.660769828692342
$ws.Range("K7").Value = 15.42909643670925
$ws.Range("L7").Value = 9.271719786232428
$ws.Range("N7").Value = 17.54899473112698
$ws.Range("O7").Value = 25.89615151554627

# Row 8
$ws.Range("C8").Value = 13.46847482466085
$ws.Range("D8").Value = 5.669170405901534
$ws.Range("E8").Value = 13.85745262397417
$ws.Range("F8").Value = 28.94714946029437
$ws.Range("G8").Value = 3.657444849575098
$ws.Range("K8").Value = 16.0084006604178
$ws.Range("L8").Value = 9.265487864698271
$ws.Range("N8").Value = 17.48830726877335
$ws.Range("O8").Value = 25.82318537872062

# Row 9
$ws.Range("C9").Value = 13.72366370511981
$ws.Range("D9").Value = 5.689812545538047
$ws.Range("E9").Value = 14.01946960210077
$ws.Range("F9").Value = 29.02288467992739
$ws.Range("G9").Value = 3.651571742874673
$ws.Range("K9").Value = 17.10205645014181
$ws.Range("L9").Value = 9.265601174507337
$ws.Range("N9").Value = 17.38427867137146
$ws.Range("O9").Value = 25.73580267642385

# Row 10
$ws.Range("C10").Value = 13.924287890285
$ws.Range("D10").Value = 5.705188832108661
$ws.Range("E10").Value = 14.15284332746969
$ws.Range("F10").Value = 29.12177037061203
$ws.Range("G10").Value = 3.647646881822344
$ws.Range("K10").Value = 17.87123050880141
$ws.Range("L10").Value = 9.273178019919971
$ws.Range("N10").Value = 17.31693111308071
$ws.Range("O10").Value = 25.70588813583062

# Row 11
$ws.Range("C11").Value = 14.01809134540233
$ws.Range("D11").Value = 5.712222504178399
$ws.Range("E11").Value = 14.21646332253239
$ws.Range("F11").Value = 29.1760902204844
$ws.Range("G11").Value = 3.645945128302705
$ws.Range("K11").Value = 18.21212497195892
$ws.Range("L11").Value = 9.278243853231302
$ws.Range("N11").Value = 17.28825502815365
$ws.Range("O11").Value = 25.69977976975739

# Row 12
$ws.Range("C12").Value = 14.05394937824802
$ws.Range("D12").Value = 5.714890940775568
$ws.Range("E12").Value = 14.24096287261803
$ws.Range("F12").Value = 29.19799423871182
$ws.Range("G12").Value = 3.645312680760316
$ws.Range("K12").Value = 18.33980260067799
$ws.Range("L12").Value = 9.280394003658303
$ws.Range("N12").Value = 17.27767735639168
$ws.Range("O12").Value = 25.6985484684655

# Row 13
$ws.Range("C13").Value = 14.04621219152628
$ws.Range("D13").Value = 5.714316036128453
$ws.Range("E13").Value = 14.23566856839892
$ws.Range("F13").Value = 29.19321762754586
$ws.Range("G13").Value = 3.645448358322838
$ws.Range("K13").Value = 18.31236946242504
$ws.Range("L13").Value = 9.279920637014966
$ws.Range("N13").Value = 17.27994294391954
$ws.Range("O13").Value = 25.69876549308464

# Row 14
$ws.Range("C14").Value = 14.02103482366224
$ws.Range("D14").Value = 5.712441942892081
$ws.Range("E14").Value = 14.21847083085117
$ws.Range("F14").Value = 29.17786559102387
$ws.Range("G14").Value = 3.645892856940917
$ws.Range("K14").Value = 18.22265800801623
$ws.Range("L14").Value = 9.278416107412447
$ws.Range("N14").Value = 17.28737916162664
$ws.Range("O14").Value = 25.69965677263351

# Row 15
$ws.Range("C15").Value = 14.00565595653415
$ws.Range("D15").Value = 5.711294630083982
$ws.Range("E15").Value = 14.207989376267
$ws.Range("F15").Value = 29.16863550684098
$ws.Range("G15").Value = 3.646166682071583
$ws.Range("K15").Value = 18.16752003731848
$ws.Range("L15").Value = 9.277524698597668
$ws.Range("N15").Value = 17.29197068096161
$ws.Range("O15").Value = 25.70034367526362

# Row 16
$ws.Range("C16").Value = 13.91820628664606
$ws.Range("D16").Value = 5.70472989200896
$ws.Range("E16").Value = 14.14874358936847
$ws.Range("F16").Value = 29.11840764893252
$ws.Range("G16").Value = 3.647759773877465
$ws.Range("K16").Value = 17.84876079682871
$ws.Range("L16").Value = 9.272879454428832
$ws.Range("N16").Value = 17.31884456416005
$ws.Range("O16").Value = 25.70643855150437

# Row 17
$ws.Range("C17").Value = 13.86518895484645
$ws.Range("D17").Value = 5.70071213848678
$ws.Range("E17").Value = 14.11314215866253
$ws.Range("F17").Value = 29.0899804798113
$ws.Range("G17").Value = 3.648758472075849
$ws.Range("K17").Value = 17.65081899595524
$ws.Range("L17").Value = 9.27044377486413
$ws.Range("N17").Value = 17.3358325568231
$ws.Range("O17").Value = 25.71210108204834

# Row 18
$ws.Range("C18").Value = 13.83493546355751
$ws.Range("D18").Value = 5.69840493602367
$ws.Range("E18").Value = 14.09294359924058
$ws.Range("F18").Value = 29.07450932031128
$ws.Range("G18").Value = 3.649340777548811
$ws.Range("K18").Value = 17.53612659767329
$ws.Range("L18").Value = 9.269195301476504
$ws.Range("N18").Value = 17.34578818146034
$ws.Range("O18").Value = 25.71606379660227

# Row 19
$ws.Range("C19").Value = 13.824734340916
$ws.Range("D19").Value = 5.697624412984551
$ws.Range("E19").Value = 14.08615300638279
$ws.Range("F19").Value = 29.06942230772481
$ws.Range("G19").Value = 3.64953929165957
$ws.Range("K19").Value = 17.49715292790915
$ws.Range("L19").Value = 9.26879880372184
$ws.Range("N19").Value = 17.34919070553982
$ws.Range("O19").Value = 25.71752660740246

# Row 20
$ws.Range("C20").Value = 13.87080803468086
$ws.Range("D20").Value = 5.701139456705963
$ws.Range("E20").Value = 14.1169032862849
$ws.Range("F20").Value = 29.09291563973115
$ws.Range("G20").Value = 3.648651343791997
$ws.Range("K20").Value = 17.67197821149468
$ws.Range("L20").Value = 9.27068728427772
$ws.Range("N20").Value = 17.3340050570213
$ws.Range("O20").Value = 25.71142522411711

# Row 21
$ws.Range("C21").Value = 14.02842111755956
$ws.Range("D21").Value = 5.712992280764224
$ws.Range("E21").Value = 14.22351127900741
$ws.Range("F21").Value = 29.18233872047874
$ws.Range("G21").Value = 3.645761972581509
$ws.Range("K21").Value = 18.24904762374977
$ws.Range("L21").Value = 9.278851741096931
$ws.Range("N21").Value = 17.28518733342486
$ws.Range("O21").Value = 25.69936559975084

# Row 22
$ws.Range("C22").Value = 14.13337942232257
$ws.Range("D22").Value = 5.720767259184653
$ws.Range("E22").Value = 14.29555627801474
$ws.Range("F22").Value = 29.24855376970643
$ws.Range("G22").Value = 3.64394333983695
$ws.Range("K22").Value = 18.61792103918674
$ws.Range("L22").Value = 9.285538399093815
$ws.Range("N22").Value = 17.25492172857994
$ws.Range("O22").Value = 25.69779050567092

# Row 23
$ws.Range("C23").Value = 14.07719235962186
$ws.Range("D23").Value = 5.716615220111991
$ws.Range("E23").Value = 14.25689302064259
$ws.Range("F23").Value = 29.21250567372233
$ws.Range("G23").Value = 3.644907618420558
$ws.Range("K23").Value = 18.42183914914893
$ws.Range("L23").Value = 9.281846379884163
$ws.Range("N23").Value = 17.27092521999594
$ws.Range("O23").Value = 25.69805321691519

# Row 24
$ws.Range("C24").Value = 13.86826693985713
$ws.Range("D24").Value = 5.700946257784739
$ws.Range("E24").Value = 14.115202040853
$ws.Range("F24").Value = 29.09158593702509
$ws.Range("G24").Value = 3.648699751126379
$ws.Range("K24").Value = 17.66241490270146
$ws.Range("L24").Value = 9.27057672065898
$ws.Range("N24").Value = 17.33483068087679
$ws.Range("O24").Value = 25.71172857662171

# Row 25
$ws.Range("C25").Value = 13.65221614420454
$ws.Range("D25").Value = 5.684188903628498
$ws.Range("E25").Value = 13.97306688076373
$ws.Range("F25").Value = 28.99478963386195
$ws.Range("G25").Value = 3.653091748678894
$ws.Range("K25").Value = 16.81165240486593
$ws.Range("L25").Value = 9.264252154467756
$ws.Range("N25").Value = 17.41082274125345
$ws.Range("O25").Value = 25.75344205684616
